# Daily attendance processing - reorder the "Recorded By" (column G) values
# for every row where multiple recorders are listed as a comma-separated
# string: reverse the order of the comma-separated entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val.ToString().Contains(",")) {
        $text = $val.ToString()
        $parts = $text.Split(",")

        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $count = $trimmed.Length
        $reversed = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $reversed += $trimmed[$i]
        }

        $newVal = [string]::Join(", ", $reversed)
        $cell.Value2 = $newVal
    }
}
